$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, styled like the existing header row (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-19
$data = @(
    @(7, 7),
    @(9, 9),
    @(5, 6),
    @(6, 6),
    @(7, 8),
    @(6, 6),
    @(10, 10),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 9),
    @(1, 3),
    @(1, 5),
    @(4, 5),
    @(3, 5),
    @(1, 2)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
